$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns affected by the row swap: A<->row4, B<->row4, D,E,F,G,H,Q,R similarly,
# and row3<->row5 for the same columns.
$cols = @("A","B","D","E","F","G","H","Q","R")

function Swap-Rows($ws, $cols, $r1, $r2) {
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}

Swap-Rows $ws $cols 2 4
Swap-Rows $ws $cols 3 5
